$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I13").Value = 'aa'
$ws.Range("J13").Value = 'Agree/Accept'
$ws.Range("I22").Value = '%'
$ws.Range("J22").Value = 'Uninterpretable'
$ws.Range("I33").Value = 'sd'
$ws.Range("J33").Value = 'Statement-non-opinion'
$ws.Range("I34").Value = 'sd'
$ws.Range("J34").Value = 'Statement-non-opinion'
$ws.Range("I38").Value = 'aa'
$ws.Range("J38").Value = 'Agree/Accept'
$ws.Range("I51").Value = 'aa'
$ws.Range("J51").Value = 'Agree/Accept'
$ws.Range("I52").Value = '%'
$ws.Range("J52").Value = 'Uninterpretable'
$ws.Range("I73").Value = 'sd'
$ws.Range("J73").Value = 'Statement-non-opinion'
$ws.Range("I95").Value = 'sv'
$ws.Range("J95").Value = 'Statement-opinion'
$ws.Range("I104").Value = 'sd'
$ws.Range("J104").Value = 'Statement-non-opinion'
$ws.Range("I111").Value = 'aa'
$ws.Range("J111").Value = 'Agree/Accept'
$ws.Range("I112").Value = 'aa'
$ws.Range("J112").Value = 'Agree/Accept'
$ws.Range("I130").Value = 'sv'
$ws.Range("J130").Value = 'Statement-opinion'
$ws.Range("I131").Value = 'sd'
$ws.Range("J131").Value = 'Statement-non-opinion'
$ws.Range("I146").Value = 'sd'
$ws.Range("J146").Value = 'Statement-non-opinion'
$ws.Range("I150").Value = 'b'
$ws.Range("J150").Value = 'Acknowledge (Backchannel)'
$ws.Range("I152").Value = 'aa'
$ws.Range("J152").Value = 'Agree/Accept'
$ws.Range("I163").Value = 'aa'
$ws.Range("J163").Value = 'Agree/Accept'
$ws.Range("I164").Value = 'aa'
$ws.Range("J164").Value = 'Agree/Accept'
$ws.Range("I182").Value = 'sd'
$ws.Range("J182").Value = 'Statement-non-opinion'
$ws.Range("I189").Value = 'ba'
$ws.Range("J189").Value = 'Appreciation'
$ws.Range("I207").Value = 'aa'
$ws.Range("J207").Value = 'Agree/Accept'
$ws.Range("I211").Value = 'sv'
$ws.Range("J211").Value = 'Statement-opinion'
$ws.Range("I216").Value = 'b'
$ws.Range("J216").Value = 'Acknowledge (Backchannel)'
$ws.Range("I223").Value = 'aa'
$ws.Range("J223").Value = 'Agree/Accept'
$ws.Range("I224").Value = 'aa'
$ws.Range("J224").Value = 'Agree/Accept'
$ws.Range("I227").Value = 'sd'
$ws.Range("J227").Value = 'Statement-non-opinion'
$ws.Range("I247").Value = 'aa'
$ws.Range("J247").Value = 'Agree/Accept'
$ws.Range("I266").Value = 'aa'
$ws.Range("J266").Value = 'Agree/Accept'
$ws.Range("I267").Value = 'sd'
$ws.Range("J267").Value = 'Statement-non-opinion'
$ws.Range("I279").Value = 'ba'
$ws.Range("J279").Value = 'Appreciation'
$ws.Range("I320").Value = 'aa'
$ws.Range("J320").Value = 'Agree/Accept'
$ws.Range("I330").Value = 'qy'
$ws.Range("J330").Value = 'Yes-No-Question'
$ws.Range("I343").Value = '%'
$ws.Range("J343").Value = 'Uninterpretable'
$ws.Range("I352").Value = 'sd'
$ws.Range("J352").Value = 'Statement-non-opinion'
$ws.Range("I353").Value = 'ba'
$ws.Range("J353").Value = 'Appreciation'
$ws.Range("I354").Value = 'aa'
$ws.Range("J354").Value = 'Agree/Accept'
$ws.Range("I355").Value = 'sd'
$ws.Range("J355").Value = 'Statement-non-opinion'
$ws.Range("I359").Value = 'sv'
$ws.Range("J359").Value = 'Statement-opinion'
$ws.Range("I376").Value = 'sd'
$ws.Range("J376").Value = 'Statement-non-opinion'
$ws.Range("I381").Value = 'sv'
$ws.Range("J381").Value = 'Statement-opinion'
$ws.Range("I384").Value = 'sd'
$ws.Range("J384").Value = 'Statement-non-opinion'
$ws.Range("I389").Value = 'sv'
$ws.Range("J389").Value = 'Statement-opinion'
